$p = $ppt.ActivePresentation
Write-Host "Designs.Count before: " $p.Designs.Count
try {
  $newd = $p.Designs.Add("dummy")
  Write-Host "Add returned: " $newd
} catch {
  Write-Host "Add err: $_"
}
Write-Host "Designs.Count after: " $p.Designs.Count
